# Se procesan de nuevo los datos con las nuevas dimensiones curadas
# municipio-nombre (col C) pasa de ser medida a ser dimension (refArea),
# y diputados (col D) pasa de ser dimension a ser medida.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = "sdmx-dimension:refArea"
$ws.Range("D2").Value = "iaest-measure:diputados"

$ws.Range("C3").Value = "dim"
$ws.Range("D3").Value = "medida"

$ws.Range("C4").Value = "URI-Municipio"
$ws.Range("D4").Value = "xsd:int"

$ws.Range("D5").Value = ""
